$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Original data rows (2-5):
#   2: even_MAG-GUT1861.fa
#   3: even_MAG-GUT43440.fa
#   4: even_MAG-GUT91702.fa
#   5: even_MAG-GUT91898.fa
#
# Target data rows (2-3):
#   2: even_MAG-GUT43440.fa   (was row 3)
#   3: even_MAG-GUT91898.fa   (was row 5)
#
# So rows 2 and 4 (even_MAG-GUT1861.fa and even_MAG-GUT91702.fa) must be
# removed, and the remaining rows shift up to close the gaps.

# Delete row 2 (even_MAG-GUT1861.fa). Remaining rows shift up:
#   old row3 -> row2, old row4 -> row3, old row5 -> row4
$ws.Rows.Item(2).Delete()

# After the shift above, even_MAG-GUT91702.fa is now row 3. Delete it.
# Remaining rows shift up again: old row5 -> row3
$ws.Rows.Item(3).Delete()
